# This edit inserts one new data row into Sheet1 at row 370, which pushes
# the existing rows 370-485 down to 371-486 (a weekly price-update style
# edit - "Fruta / hortaliza, semanal"). The new row carries a fresh price
# record for "Brocoli" at "Feria Lagunitas de Puerto Montt".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 370; Excel automatically shifts rows
# 370:485 down to 371:486 and copies the row formatting (e.g. the date
# style used in column D) onto the new row.
$ws.Rows.Item(370).Insert()

# Populate the newly inserted row 370 with the new record's values.
$ws.Cells.Item(370, 1).Value2  = 4
$ws.Cells.Item(370, 2).Value2  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(370, 3).Value2  = 'Los Lagos'
$ws.Cells.Item(370, 4).Value2  = 44985
$ws.Cells.Item(370, 5).Value2  = 10
$ws.Cells.Item(370, 6).Value2  = 100112023
$ws.Cells.Item(370, 7).Value2  = 'Brócoli'
$ws.Cells.Item(370, 8).Value2  = 'Sin especificar'
$ws.Cells.Item(370, 9).Value2  = 'Primera'
$ws.Cells.Item(370, 10).Value2 = 1000
$ws.Cells.Item(370, 11).Value2 = 1600
$ws.Cells.Item(370, 12).Value2 = 1600
$ws.Cells.Item(370, 13).Value2 = 1600
$ws.Cells.Item(370, 14).Value2 = '$/unidad'
$ws.Cells.Item(370, 15).Value2 = 'Región Metropolitana'
$ws.Cells.Item(370, 16).Value2 = 1600
$ws.Cells.Item(370, 17).Value2 = 1
$ws.Cells.Item(370, 18).Value2 = 'Hortaliza'
